$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$firstRow = $used.Row

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
